# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# per commit 'Updated cryptos list ... with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.648.79'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '1.632.65'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.95'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("E6").Value = '  +3.09%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.14'
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("D12").Value = '1.859.94'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = '1.614.67'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.09'
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = '26.652.20'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.23'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.05'
$ws.Range("E19").Value = '  +8.05%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  +2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.38'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.92'
$ws.Range("E24").Value = '  +2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.14'
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.88'
$ws.Range("E28").Value = '  +4.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.47'
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0505'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +3.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.98'
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '1.208.40'
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.807'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.28'
$ws.Range("E41").Value = '  -1.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.41'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").Value = '1.772.58'
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.56'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.56'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.72'
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0514'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.59'
$ws.Range("E49").Value = '  +4.32%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +0.14%  '
